$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'capri mens joggers'
$ws.Cells.Item(2, 1).Value = 'easton baseball pants mens'
$ws.Cells.Item(3, 1).Value = 'mueller knee pads basketball'
$ws.Cells.Item(4, 1).Value = 'bmx knee pads youth'
$ws.Cells.Item(5, 1).Value = 'basketball tights kids'
$ws.Cells.Item(6, 1).Value = 'tortoise knee pads'
$ws.Cells.Item(7, 1).Value = 'bendable knee pads'
$ws.Cells.Item(8, 1).Value = 'dakine knee pads'
$ws.Cells.Item(9, 1).Value = '187 knee pads'
$ws.Cells.Item(10, 1).Value = 'mma knee pads'
$ws.Cells.Item(11, 1).Value = 'scrubs men pants'
$ws.Cells.Item(12, 1).Value = 'copper compression pants'
$ws.Cells.Item(13, 1).Value = 'elbow pads knee pads'
$ws.Cells.Item(14, 1).Value = 'knee pad for scooter'
$ws.Cells.Item(15, 1).Value = 'hamstring compression pants'
$ws.Cells.Item(16, 1).Value = 'compression tights men basketball'
$ws.Cells.Item(17, 1).Value = 'white leggings for men'
$ws.Cells.Item(18, 1).Value = 'red leggings men'
$ws.Cells.Item(19, 1).Value = 'purple leggings men'
$ws.Cells.Item(20, 1).Value = 'yellow leggings men'
$ws.Cells.Item(21, 1).Value = 'athletic leggings for women'
$ws.Cells.Item(22, 1).Value = 'athletic leggings kids'
$ws.Cells.Item(23, 1).Value = 'basketball leggings for women'
$ws.Cells.Item(24, 1).Value = 'basketball leggings kids'
$ws.Cells.Item(25, 1).Value = 'bdu pants with knee pads'
$ws.Cells.Item(26, 1).Value = 'nike leggings for men'
$ws.Cells.Item(27, 1).Value = 'compression pants basketball'
$ws.Cells.Item(28, 1).Value = 'compression pants for women'
$ws.Cells.Item(29, 1).Value = 'compression pants kids'
$ws.Cells.Item(30, 1).Value = 'compression pants knee'
$ws.Cells.Item(31, 1).Value = 'compression pants men under armour'
$ws.Cells.Item(32, 1).Value = 'compression pants set'
$ws.Cells.Item(33, 1).Value = 'compression pants tesla'
$ws.Cells.Item(34, 1).Value = 'camo pants with knee pads'
$ws.Cells.Item(35, 1).Value = 'kids basketball pads'
$ws.Cells.Item(36, 1).Value = 'gold leggings men'
$ws.Cells.Item(37, 1).Value = 'kickboxing knee pads'
$ws.Cells.Item(38, 1).Value = 'men gym pants'
$ws.Cells.Item(39, 1).Value = 'men nike compression pants'
$ws.Cells.Item(40, 1).Value = 'swim leggings for men'
$ws.Cells.Item(41, 1).Value = 'fleece leggings men'
$ws.Cells.Item(42, 1).Value = 'sliding shorts with knee pads'
$ws.Cells.Item(43, 1).Value = 'mens basketball joggers'
$ws.Cells.Item(44, 1).Value = 'mcdavid basketball knee'
$ws.Cells.Item(45, 1).Value = '3xl knee pads'
$ws.Cells.Item(46, 1).Value = '511 knee pads'
$ws.Cells.Item(47, 1).Value = 'caving knee pads'
$ws.Cells.Item(48, 1).Value = 'enduro knee pads'
$ws.Cells.Item(49, 1).Value = 'bodyprox knee pads'
$ws.Cells.Item(50, 1).Value = 'blackhawk knee pads'
$ws.Cells.Item(51, 1).Value = 'snickers knee pads'
$ws.Cells.Item(52, 1).Value = 'bball knee pads'
$ws.Cells.Item(53, 1).Value = '661 knee pads'
$ws.Cells.Item(54, 1).Value = 'basket knee pads'
$ws.Cells.Item(55, 1).Value = 'armadillo knee pads'
$ws.Cells.Item(56, 1).Value = 'swim pants men'
$ws.Cells.Item(57, 1).Value = 'athletic capris'
$ws.Cells.Item(58, 1).Value = 'bunheads knee pads'
$ws.Cells.Item(59, 1).Value = '6xl compression pants'
$ws.Cells.Item(60, 1).Value = 'kali knee pads'
$ws.Cells.Item(61, 1).Value = 'arcteryx knee pads'
$ws.Cells.Item(62, 1).Value = 'spelunking knee pads'
$ws.Cells.Item(63, 1).Value = 'bcg compression pants'
$ws.Cells.Item(64, 1).Value = 'bcg knee pads'
$ws.Cells.Item(65, 1).Value = 'tesla mens leggings'
$ws.Cells.Item(66, 1).Value = 'awp knee pads'
$ws.Cells.Item(67, 1).Value = 'leatt knee pads'
$ws.Cells.Item(68, 1).Value = 'fr knee pads'
$ws.Cells.Item(69, 1).Value = 'eurotard knee pads'
$ws.Cells.Item(70, 1).Value = 'asics leggings men'
$ws.Cells.Item(71, 1).Value = 'alleson baseball pants youth'
$ws.Cells.Item(72, 1).Value = 'frozen knee pads'
$ws.Cells.Item(73, 1).Value = 'training tights men'
$ws.Cells.Item(74, 1).Value = 'training leggings for men'
$ws.Cells.Item(75, 1).Value = 'nba knee pads'
$ws.Cells.Item(76, 1).Value = 'prayer knee pads'
$ws.Cells.Item(77, 1).Value = 'neoprene knee pad'
$ws.Cells.Item(78, 1).Value = 's1 knee pads'
$ws.Cells.Item(79, 1).Value = 'white basketball pants'
$ws.Cells.Item(80, 1).Value = 'lululemon compression pants'
$ws.Cells.Item(81, 1).Value = 'xlarge knee pads'
$ws.Cells.Item(82, 1).Value = 'ballet tights for men'
$ws.Cells.Item(83, 1).Value = 'nike compression pants for men'
$ws.Cells.Item(84, 1).Value = 'mens compression pants adidas'
$ws.Cells.Item(85, 1).Value = 'mens compression pants champion'
$ws.Cells.Item(86, 1).Value = 'baseball pants mens with piping'
$ws.Cells.Item(87, 1).Value = 'exercise compression leggings'
$ws.Cells.Item(88, 1).Value = 'teen knee pads'
$ws.Cells.Item(89, 1).Value = 'short tights for men'
$ws.Cells.Item(90, 1).Value = 'runner compression pants'
$ws.Cells.Item(91, 1).Value = 'swimming tights men'
$ws.Cells.Item(92, 1).Value = 'half tights men'
$ws.Cells.Item(93, 1).Value = 'fishing knee pads'
$ws.Cells.Item(94, 1).Value = 'marvel knee pads'
$ws.Cells.Item(95, 1).Value = 'bauer compression pants'
$ws.Cells.Item(96, 1).Value = 'police knee pads'
$ws.Cells.Item(97, 1).Value = 'spiderman tights men'
$ws.Cells.Item(98, 1).Value = 'wwe knee pads'
$ws.Cells.Item(99, 1).Value = 'petite compression leggings'
$ws.Cells.Item(100, 1).Value = '3x compression leggings'
